$wb = $excel.ActiveWorkbook

# Column F ("想去人数" / want-to-go count) updates, applied to both the
# "展览" and "全部类型" worksheets (which carry the same data table).
$updates = @{
    2  = 149
    3  = 1693
    4  = 782
    7  = 11865
    10 = 473
    11 = 403
    14 = 13447
    15 = 13370
    20 = 272
    23 = 157
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
